$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the new "date_time" column (C) ---
# Data rows first (C2:C101) so that the "plain" date-format style is
# interned before the "right-aligned" header style, matching the
# cellXfs order produced by the source workbook (numFmt 166 plain, then
# numFmt 166 + right alignment).
$dateSerial = 43831

for ($row = 2; $row -le 101; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    $cell.Value = $dateSerial
    $cell.NumberFormat = "dd.mm.yyyy"
}

# Header cell: same date number format, but bold-ish/right aligned like
# the rest of the header styling used by this workbook.
$header = $ws.Range("C1")
$header.Value = "date_time"
$header.NumberFormat = "dd.mm.yyyy"
$header.HorizontalAlignment = -4152

# Match column A/B's explicit width.
$ws.Columns("C").ColumnWidth = 15.15

# The first two data rows shrink to the same row height already used by
# every other data row once they pick up real content in the new column.
$ws.Rows("2:3").RowHeight = 13.8

# Restore the view to the top of the sheet and park the selection where
# the edit left it.
$excel.Goto($ws.Range("A1"), $true) | Out-Null
$ws.Range("H104").Select() | Out-Null
